$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UniqueValues")

$sortRange = $ws.Range("A2:E56")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B56"), 0, 1, $null, 0)
$ws.Sort.SortFields.Add($ws.Range("A2:A56"), 0, 1, $null, 0)

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

$ws.Rows.Item(40).RowHeight = 15
$ws.Rows.Item(55).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 67.5
$ws.Rows.Item(54).RowHeight = 45
